# Insert a new weekly price-report row for "Poroto verde" (Feria Lagunitas de
# Puerto Montt) ahead of the existing row 79. Excel shifts rows 79-172 down
# to 80-173 (carrying formatting, incl. the date style on column D), and the
# new row 79 is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by inserting a brand-new row at 79.
$ws.Rows(79).Insert()

# Populate the newly inserted row 79 with the new record.
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 45195
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112031
$ws.Range("G79").Value = "Poroto verde"
$ws.Range("H79").Value = "Magnum"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 45
$ws.Range("K79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("M79").Value = 35000
$ws.Range("N79").Value = "$/malla 25 kilos"
$ws.Range("O79").Value = "Perú"
$ws.Range("P79").Value = 1400
$ws.Range("Q79").Value = 25
$ws.Range("R79").Value = "Hortaliza"
